$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    4  = @{ C = 4.317520552235576;  E = 3.856664261949816 }
    5  = @{ C = 10.48563750975209;  E = 5.658864198748459 }
    6  = @{ C = 4.65880603412161;   E = 5.409232631930561 }
    7  = @{ C = -0.8693696108860949; E = 2.311254283099773 }
    8  = @{ C = 4.057580120451165;  E = 3.242274909585041 }
    9  = @{ C = 3.942841799197594;  E = 3.826653192455631 }
    10 = @{ C = 2.435545128806416;  E = 4.047617271894799 }
    11 = @{ C = 4.218672000695523;  E = 4.02808863798465 }
    12 = @{ C = 4.174017263680696;  E = 2.074691389445271 }
    13 = @{ C = 1.699348375745302;  E = 3.362240252406901 }
    14 = @{ C = -2.856524424985296; E = 0.1239622353166103 }
    15 = @{ C = 6.240787792289715;  E = 4.112367048326182 }
    16 = @{ C = 3.815916106066686;  E = 2.762663830671319 }
    17 = @{ C = 0.757583445265464;  E = 2.564429185896056 }
    18 = @{ C = -0.152046383567539; E = 2.659925441240518 }
    19 = @{ C = -1.93082584212636;  E = 2.2373336846083 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
